# Apply updated view/order counts (column F) across sheets, mirroring the
# commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F11").Value = 379
$wsExhibit.Range("F13").Value = 38
$wsExhibit.Range("F15").Value = 106
$wsExhibit.Range("F18").Value = 5806
$wsExhibit.Range("F24").Value = 184

# Sheet "演出" (Performances)
$wsPerform = $wb.Worksheets.Item("演出")
$wsPerform.Range("F7").Value = 497
$wsPerform.Range("F14").Value = 48

# Sheet "全部类型" (All types - aggregated)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F12").Value = 497
$wsAll.Range("F13").Value = 497
$wsAll.Range("F24").Value = 379
$wsAll.Range("F27").Value = 38
$wsAll.Range("F30").Value = 106
$wsAll.Range("F34").Value = 5806
$wsAll.Range("F35").Value = 48
$wsAll.Range("F42").Value = 184
